# Update crypto price/volume data per latest scrape (GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.460.47"
$ws.Range("E2").Value = "  -0.64%  "

$ws.Range("D3").Value = "'2.096.01"
$ws.Range("E3").Value = "  -0.51%  "

$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.24%  "

$ws.Range("D5").Value = "'329.63"
$ws.Range("E5").Value = "  -0.12%  "

$ws.Range("E6").Value = "  -0.09%  "

$ws.Range("D7").Value = "'0.5209"

$ws.Range("D8").Value = "'0.4347"
$ws.Range("E8").Value = "  -0.78%  "

$ws.Range("D9").Value = "'53.16"
$ws.Range("E9").Value = "  +17.95%  "

$ws.Range("D10").Value = "'0.08858"
$ws.Range("E10").Value = "  -1.19%  "

$ws.Range("E11").Value = "  -1.54%  "

$ws.Range("D12").Value = "'24.37"
$ws.Range("E12").Value = "  -1.78%  "

$ws.Range("D13").Value = "'2.089.06"
$ws.Range("E13").Value = "  -0.90%  "

$ws.Range("D14").Value = "'6.682"
$ws.Range("E14").Value = "  -1.50%  "

$ws.Range("D15").Value = "'7.673"
$ws.Range("E15").Value = "  +0.29%  "

$ws.Range("D16").Value = "'95.76"

$ws.Range("E17").Value = "  +0.00%  "

$ws.Range("D18").Value = "'0.00001119"
$ws.Range("E18").Value = "  -1.41%  "

$ws.Range("D19").Value = "'0.06581"
$ws.Range("E19").Value = "  -0.31%  "

$ws.Range("D20").Value = "'19.21"
$ws.Range("E20").Value = "  +0.32%  "

$ws.Range("E21").Value = "  -0.13%  "

$ws.Range("E22").Value = "  -2.29%  "

$ws.Range("D23").Value = "'30.491.61"
$ws.Range("E23").Value = "  -1.12%  "

$ws.Range("D24").Value = "'12.19"
$ws.Range("E24").Value = "  +2.07%  "

$ws.Range("D25").Value = "'2.338"
$ws.Range("E25").Value = "  +3.29%  "

$ws.Range("D26").Value = "'2.334.56"
$ws.Range("E26").Value = "  -1.10%  "

$ws.Range("D27").Value = "'22.26"
$ws.Range("E27").Value = "  -2.57%  "

$ws.Range("D28").Value = "'2.582"
$ws.Range("E28").Value = "  +1.94%  "

$ws.Range("D29").Value = "'162.38"
$ws.Range("E29").Value = "  -0.72%  "

$ws.Range("D30").Value = "'131.62"
$ws.Range("E30").Value = "  -1.53%  "

$ws.Range("D31").Value = "'1.190"
$ws.Range("E31").Value = "  +0.48%  "

$ws.Range("E32").Value = "  -0.06%  "

$ws.Range("D33").Value = "'1.675"
$ws.Range("E33").Value = "  +10.44%  "

$ws.Range("D34").Value = "'6.132"
$ws.Range("E34").Value = "  -1.12%  "

$ws.Range("E35").Value = "  -0.87%  "

$ws.Range("E36").Value = "  +5.98%  "

$ws.Range("D37").Value = "'0.02573"
$ws.Range("E37").Value = "  -0.51%  "

$ws.Range("D38").Value = "'0.06807"
$ws.Range("E38").Value = "  +0.68%  "

$ws.Range("D39").Value = "'12.75"
$ws.Range("E39").Value = "  -0.01%  "

$ws.Range("D40").Value = "'5.457"
$ws.Range("E40").Value = "  -2.57%  "

$ws.Range("E41").Value = "  +0.77%  "

$ws.Range("D42").Value = "'0.6906"
$ws.Range("E42").Value = "  +2.08%  "

$ws.Range("E43").Value = "  +1.12%  "

$ws.Range("D44").Value = "'1.000"
$ws.Range("E44").Value = "  -0.10%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'13.96"
$ws.Range("E45").Value = "  -1.41%  "

$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.6366"
$ws.Range("E46").Value = "  +1.20%  "

$ws.Range("E47").Value = "  -2.12%  "

$ws.Range("D48").Value = "'3.622"
$ws.Range("E48").Value = "  -0.71%  "

$ws.Range("D49").Value = "'1.238"
$ws.Range("E49").Value = "  +7.55%  "

$ws.Range("D50").Value = "'1.240"
$ws.Range("E50").Value = "  -2.47%  "

$ws.Range("D51").Value = "'81.87"
$ws.Range("E51").Value = "  -1.23%  "
